$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the slurm job id notes for the two existing rows (55, 56)
$ws.Range("P55").Value2 = "slurm-42599447"
$ws.Range("P56").Value2 = "slurm-42599773"

# Append a new row 57, mirroring row 56 but for "16 words" (column B)
$ws.Range("A57").Value2 = $ws.Range("A56").Value2
$ws.Range("B57").Value2 = "16 words"
$ws.Range("C57").Value2 = $ws.Range("C56").Value2
$ws.Range("D57").Value2 = $ws.Range("D56").Value2
$ws.Range("E57").Value2 = $ws.Range("E56").Value2
$ws.Range("F57").Value2 = $ws.Range("F56").Value2
$ws.Range("G57").Value2 = $ws.Range("G56").Value2
$ws.Range("H57").Value2 = $ws.Range("H56").Value2
$ws.Range("I57").Value2 = $ws.Range("I56").Value2
$ws.Range("J57").Value2 = $ws.Range("J56").Value2
$ws.Range("K57").Value2 = $ws.Range("K56").Value2
$ws.Range("P57").Value2 = "slurm-42599776"

# Update the saved view state to match the new selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 36
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("H44").Select()
